$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Prefers sustainable future"
$ws.Range("B2").Value = 0.680881448179833
$ws.Range("D2").Value = 0.721220700610068
$ws.Range("K2").Value = 0.759539408850822
$ws.Range("L2").Value = 0.690381293922693
$ws.Range("N2").Value = 0.618731104894089

# Row 3 - "Prefers sustainable future (Variant: Scenario A = Sustainable)"
$ws.Range("B3").Value = 0.69861761980374
$ws.Range("K3").Value = 0.757178577397713
$ws.Range("L3").Value = 0.733729575456935
$ws.Range("N3").Value = 0.623193310926974

# Row 4 - "Prefers sustainable future (Variant: Scenario B = Sustainable)"
$ws.Range("B4").Value = 0.662577293517367
$ws.Range("K4").Value = 0.759030706028559
$ws.Range("L4").Value = 0.645983607142547
$ws.Range("N4").Value = 0.610255864577834
